# Weekly fruit/vegetable price update: insert a new weekly record.
# A new row is inserted right after the existing row 17 (which holds the
# most recent-looking record at the time), pushing all subsequent rows
# down by one. The new row duplicates row 17's price data but carries a
# new date (2022-08-30, serial 44803).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 18; rows 18-39 shift down to 19-40.
$ws.Rows("18:18").Insert()

# Populate the newly inserted row 18 with the new weekly entry.
$ws.Range("A18").Value = 4
$ws.Range("B18").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C18").Value = "Los Lagos"
$ws.Range("D18").Value = 44803
$ws.Range("E18").Value = 10
$ws.Range("F18").Value = 100112012
$ws.Range("G18").Value = "Espinaca"
$ws.Range("H18").Value = "Sin especificar"
$ws.Range("I18").Value = "Primera"
$ws.Range("J18").Value = 35
$ws.Range("K18").Value = 12000
$ws.Range("L18").Value = 12000
$ws.Range("M18").Value = 12000
$ws.Range("N18").Value = "$/cuna 10 kilos"
$ws.Range("O18").Value = "Región Metropolitana"
$ws.Range("P18").Value = 1200
$ws.Range("Q18").Value = 10
$ws.Range("R18").Value = "Hortaliza"
